$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 84618950
$ws.Range("I70").Value = 16668481
$ws.Range("K70").Value = 50005443
$ws.Range("M70").Value = -50005173
$ws.Range("H73").Value = 84618950
$ws.Range("I73").Value = 16668481
$ws.Range("K73").Value = 50005443
$ws.Range("M73").Value = -50004507
$ws.Range("H131").Value = 9647.583000000001
$ws.Range("I131").Value = 5063
$ws.Range("K131").Value = 15189
$ws.Range("M131").Value = -10149
$ws.Range("H137").Value = 4499.5
$ws.Range("I137").Value = 4999
$ws.Range("K137").Value = 14997
$ws.Range("M137").Value = -12447
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 131519.36
$ws.Range("I32").Value = 131519.36
$ws.Range("K32").Value = 131519.36
$ws.Range("M32").Value = -131232.36
$ws.Range("H45").Value = 4933.6216
$ws.Range("I45").Value = 5848.08
$ws.Range("K45").Value = 5848.08
$ws.Range("M45").Value = -5471.08
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 574.25
$ws.Range("I64").Value = 696.1667
$ws.Range("J64").Value = 208.5
$ws.Range("K64").Value = 696.1667
$ws.Range("L64").Value = 208.5
$ws.Range("M64").Value = -471.1667
$ws.Range("N64").Value = -658.5
$ws.Range("H67").Value = 574.25
$ws.Range("I67").Value = 696.1667
$ws.Range("J67").Value = 208.5
$ws.Range("K67").Value = 696.1667
$ws.Range("L67").Value = 208.5
$ws.Range("M67").Value = 83.83330000000001
$ws.Range("N67").Value = -1768.5
$ws.Range("H134").Value = 2601.1035
$ws.Range("I134").Value = 2237.28
$ws.Range("K134").Value = 6711.84
$ws.Range("M134").Value = -4176.84
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1376.1578
$ws.Range("I16").Value = 1282.0714
$ws.Range("J16").Value = 1639.6
$ws.Range("K16").Value = 1282.0714
$ws.Range("L16").Value = 1639.6
$ws.Range("M16").Value = -995.0714
$ws.Range("N16").Value = -2213.6
$ws.Range("H58").Value = 2276.8215
$ws.Range("I58").Value = 2141.762
$ws.Range("K58").Value = 2141.762
$ws.Range("M58").Value = -1938.762
$ws.Range("H113").Value = 1376.1578
$ws.Range("I113").Value = 1282.0714
$ws.Range("J113").Value = 1639.6
$ws.Range("K113").Value = 1282.0714
$ws.Range("L113").Value = 1639.6
$ws.Range("M113").Value = 887.9286
$ws.Range("N113").Value = -5979.6
$ws.Range("H136").Value = 2276.8215
$ws.Range("I136").Value = 2141.762
$ws.Range("K136").Value = 6425.286
$ws.Range("M136").Value = -3875.286
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 18
$ws.Range("I2").Value = 18
$ws.Range("K2").Value = 108
$ws.Range("M2").Value = 5
$ws.Range("H94").Value = 8825.125
$ws.Range("I94").Value = 1000.5
$ws.Range("J94").Value = 11433.333
$ws.Range("K94").Value = 3001.5
$ws.Range("L94").Value = 34299.999
$ws.Range("M94").Value = -2325.5
$ws.Range("N94").Value = -35651.999
$ws.Range("H132").Value = 2129.6897
$ws.Range("J132").Value = 2713.9285
$ws.Range("L132").Value = 24425.3565
$ws.Range("N132").Value = -29485.3565
$ws.Range("H133").Value = 11553.8
$ws.Range("I133").Value = 5271.5
$ws.Range("K133").Value = 15814.5
$ws.Range("M133").Value = -10754.5
$ws.Range("H137").Value = 5399.2
$ws.Range("I137").Value = 5571.7144
$ws.Range("J137").Value = 4996.6665
$ws.Range("K137").Value = 16715.1432
$ws.Range("L137").Value = 14989.9995
$ws.Range("M137").Value = -11615.1432
$ws.Range("N137").Value = -25189.9995
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 1000000000
$ws.Range("J64").Value = 1000000000
$ws.Range("L64").Value = 1000000000
$ws.Range("N64").Value = -1000000496
$ws.Range("H67").Value = 1000000000
$ws.Range("J67").Value = 1000000000
$ws.Range("L67").Value = 1000000000
$ws.Range("N67").Value = -1000001716
$ws.Range("H102").Value = 26319632
$ws.Range("I102").Value = 45456640
$ws.Range("K102").Value = 45456640
$ws.Range("M102").Value = -45455018
$ws.Range("H107").Value = 1935.88
$ws.Range("I107").Value = 879.125
$ws.Range("J107").Value = 2433.1765
$ws.Range("K107").Value = 879.125
$ws.Range("L107").Value = 2433.1765
$ws.Range("M107").Value = 1040.875
$ws.Range("N107").Value = -6273.1765
$ws.Range("H122").Value = 6899.864
$ws.Range("I122").Value = 5267.3076
$ws.Range("K122").Value = 15801.9228
$ws.Range("M122").Value = -13351.9228
$ws.Range("H126").Value = 3136.2
$ws.Range("I126").Value = 2545.25
$ws.Range("K126").Value = 7635.75
$ws.Range("M126").Value = -5165.75
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3999.25
$ws.Range("I7").Value = 1999
$ws.Range("K7").Value = 1999
$ws.Range("M7").Value = -1887
$ws.Range("H40").Value = 5287.4194
$ws.Range("I40").Value = 4112.647
$ws.Range("K40").Value = 4112.647
$ws.Range("M40").Value = -3976.647
$ws.Range("H68").Value = 7475.154
$ws.Range("I68").Value = 1794.25
$ws.Range("J68").Value = 10000
$ws.Range("K68").Value = 1794.25
$ws.Range("L68").Value = 10000
$ws.Range("M68").Value = -1045.25
$ws.Range("N68").Value = -11498
$ws.Range("H71").Value = 7475.154
$ws.Range("I71").Value = 1794.25
$ws.Range("J71").Value = 10000
$ws.Range("K71").Value = 8971.25
$ws.Range("L71").Value = 50000
$ws.Range("M71").Value = -5227.25
$ws.Range("N71").Value = -57488
$ws.Range("H122").Value = 10000
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H126").Value = 3999.25
$ws.Range("I126").Value = 1999
$ws.Range("K126").Value = 5997
$ws.Range("M126").Value = -3527
$ws.Range("H136").Value = 59353.453
$ws.Range("I136").Value = 15288.3
$ws.Range("J136").Value = 500005
$ws.Range("K136").Value = 45864.89999999999
$ws.Range("L136").Value = 1500015
$ws.Range("M136").Value = -43314.89999999999
$ws.Range("N136").Value = -1505115
$ws.Range("H139").Value = 83766.55499999999
$ws.Range("J139").Value = 83766.55499999999
$ws.Range("L139").Value = 83766.55499999999
$ws.Range("N139").Value = -94046.55499999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 12250.143
$ws.Range("I62").Value = 10500.667
$ws.Range("K62").Value = 10500.667
$ws.Range("M62").Value = -9876.666999999999
$ws.Range("H65").Value = 12250.143
$ws.Range("I65").Value = 10500.667
$ws.Range("K65").Value = 52503.335
$ws.Range("M65").Value = -49383.335
$ws.Range("H81").Value = 3914.1428
$ws.Range("I81").Value = 2899.8333
$ws.Range("K81").Value = 5799.6666
$ws.Range("M81").Value = -4738.6666
$ws.Range("H84").Value = 3914.1428
$ws.Range("I84").Value = 2899.8333
$ws.Range("K84").Value = 28998.333
$ws.Range("M84").Value = -23694.333
$ws.Range("H122").Value = 3198.4866
$ws.Range("I122").Value = 1222.08
$ws.Range("K122").Value = 3666.24
$ws.Range("M122").Value = -1216.24
$ws.Range("H126").Value = 2732.611
$ws.Range("I126").Value = 2775.7058
$ws.Range("K126").Value = 8327.117400000001
$ws.Range("M126").Value = -5857.117400000001
$ws.Range("H136").Value = 4554.1875
$ws.Range("I136").Value = 1699.3125
$ws.Range("K136").Value = 5097.9375
$ws.Range("M136").Value = -2547.9375
